$d = $word.ActiveDocument

$d.Content.Find.Execute('英文', $true, $false, $false, $false, $false, $true, 1, $false, 'English', 2) | Out-Null
$d.Content.Find.Execute(' / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文', $true, $false, $false, $false, $false, $true, 1, $false, ' / Portuguese / French / Thai / Vietnamese / Spanish', 2) | Out-Null
$d.Content.Find.Execute('簡介', $true, $false, $false, $false, $false, $true, 1, $false, 'Brief', 2) | Out-Null
$d.Content.Find.Execute('發送給在目標國家中已回覆參加但尚未寄送文件的合作夥伴的電子郵件。 將通過 customer.io 發送', $true, $false, $false, $false, $false, $true, 1, $false, 'An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io', 2) | Out-Null
$d.Content.Find.Execute('目標受眾', $true, $false, $false, $false, $false, $true, 1, $false, 'Target audience', 2) | Out-Null
$d.Content.Find.Execute('尚未提交文件的受邀合作夥伴', $true, $false, $false, $false, $false, $true, 1, $false, 'Invited partners who haven’t submitted their documents', 2) | Out-Null
$d.Content.Find.Execute('主題行', $true, $false, $false, $false, $false, $true, 1, $false, 'Subject line', 2) | Out-Null
$d.Content.Find.Execute('[活動名稱]', $true, $false, $false, $false, $false, $true, 1, $false, '[EVENT NAME]', 2) | Out-Null
$d.Content.Find.Execute(' — 您已提交文件了嗎？  ', $true, $false, $false, $false, $false, $true, 1, $false, ' — have you submitted your docs?  ', 2) | Out-Null
$d.Content.Find.Execute('不要忘記傳送文檔', $true, $false, $false, $false, $false, $true, 1, $false, 'Don’t forget to send your documents', 2) | Out-Null
$d.Content.Find.Execute('您好 ', $true, $false, $false, $false, $false, $true, 1, $false, 'Hi ', 2) | Out-Null
$d.Content.Find.Execute('[合作夥伴姓名]', $true, $false, $false, $false, $false, $true, 1, $false, '[PARTNER NAME]', 2) | Out-Null
$d.Content.Find.Execute('， ', $true, $false, $false, $false, $false, $true, 1, $false, ', ', 2) | Out-Null
$d.Content.Find.Execute('很高興能在即將舉行的 ', $true, $false, $false, $false, $false, $true, 1, $false, 'We’re excited to see you at the upcoming ', 2) | Out-Null
$d.Content.Find.Execute(' 見到您。 ‘', $true, $false, $false, $false, $false, $true, 1, $false, '. ‘', 2) | Out-Null
$d.Content.Find.Execute(' 見到您。 ', $true, $false, $false, $false, $false, $true, 1, $false, '. ', 2) | Out-Null
$d.Content.Find.Execute('為了確認註冊，需要您在 ', $true, $false, $false, $false, $false, $true, 1, $false, 'To confirm your registration, we need the following documents from you by ', 2) | Out-Null
$d.Content.Find.Execute('日月年', $true, $false, $false, $false, $false, $true, 1, $false, 'DD Mmm YYYY', 2) | Out-Null
$d.Content.Find.Execute(' 之前提供以下文檔：', $true, $false, $false, $false, $false, $true, 1, $false, ':', 2) | Out-Null
$d.Content.Find.Execute('[插入所需文件清單]', $true, $false, $false, $false, $false, $true, 1, $false, '[insert list of documents required]', 2) | Out-Null
$d.Content.Find.Execute('請將這些文檔的副本傳送給您的區域經理 ', $true, $false, $false, $false, $false, $true, 1, $false, 'Please send a copy of these documents to your country manager, ', 2) | Out-Null
$d.Content.Find.Execute('[姓名]', $true, $false, $false, $false, $false, $true, 1, $false, '[NAME]', 2) | Out-Null
$d.Content.Find.Execute('，郵箱地址為 ', $true, $false, $false, $false, $false, $true, 1, $false, ', at ', 2) | Out-Null
$d.Content.Find.Execute('[郵箱地址]', $true, $false, $false, $false, $false, $true, 1, $false, '[EMAIL ADDRESS]', 2) | Out-Null
$d.Content.Find.Execute(' 或 ', $true, $false, $false, $false, $false, $true, 1, $false, ' or ', 2) | Out-Null
$d.Content.Find.Execute('[WHATSAPP 號碼]', $true, $false, $false, $false, $false, $true, 1, $false, '[WHATSAPP NO]', 2) | Out-Null
$d.Content.Find.Execute(' (WhatsApp)，以便我們做出必要的安排，包括住宿和交通。', $true, $false, $false, $false, $false, $true, 1, $false, ' (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.', 2) | Out-Null
$d.Content.Find.Execute('如有任何疑問，請聯繫您的區域經理。', $true, $false, $false, $false, $false, $true, 1, $false, 'If you have any questions, please contact your country manager.', 2) | Out-Null
$d.Content.Find.Execute('期待在那裡見到您！', $true, $false, $false, $false, $false, $true, 1, $false, 'We look forward to seeing you there!', 2) | Out-Null
$d.Content.Find.Execute('尊敬的 ', $true, $false, $false, $false, $false, $true, 1, $false, 'Dear ', 2) | Out-Null
$d.Content.Find.Execute('為了確保您在此次活動中擁有最佳體驗，我們需要您在 ', $true, $false, $false, $false, $false, $true, 1, $false, 'To ensure you have the best experience at this event, we need the following documents from you by ', 2) | Out-Null
$d.Content.Find.Execute('請回覆此電子郵件，附上這些文檔的副本，以便我們做出必要的安排，包括住宿和交通。', $true, $false, $false, $false, $false, $true, 1, $false, 'Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.', 2) | Out-Null
$d.Content.Find.Execute('如有任何疑問，請通過 ', $true, $false, $false, $false, $false, $true, 1, $false, 'If you have any questions, please contact us via ', 2) | Out-Null
$d.Content.Find.Execute('即時聊天', $true, $false, $false, $false, $false, $true, 1, $false, 'live chat', 2) | Out-Null
$d.Content.Find.Execute(' 與我們聯繫。 ', $true, $false, $false, $false, $false, $true, 1, $false, '. ', 2) | Out-Null
$d.Content.Find.Execute('如有任何疑問，請聯繫您的區域經理 ', $true, $false, $false, $false, $false, $true, 1, $false, 'If you have any questions, please contact your country manager, ', 2) | Out-Null
$d.Content.Find.Execute(' (WhatsApp)。', $true, $false, $false, $false, $false, $true, 1, $false, ' (WhatsApp).', 2) | Out-Null

# Comment text update (Chinese -> English), preserving comment id/author/date
$d.Comments.Item(1).Range.Text = 'choose either one'

"done"
